$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.588675333333333
$ws.Range("H2").Value = 10.766026
$ws.Range("I2").Value = 0.1217029912931362
$ws.Range("J2").Value = 0.1217029912931362
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 84.82120407002421
$ws.Range("R2").Value = 763.390836630218
$ws.Range("S2").Value = 0.008309691103981474
$ws.Range("T2").Value = 0.008309691103981474

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.588675333333333
$ws.Range("H3").Value = 10.766026
$ws.Range("I3").Value = 0.1217029912931362
$ws.Range("J3").Value = 0.1217029912931362
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 650.5851781976138
$ws.Range("R3").Value = 5855.266603778525
$ws.Range("S3").Value = 0.06373597176465276
$ws.Range("T3").Value = 0.06373597176465275

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.588675333333333
$ws.Range("H4").Value = 10.766026
$ws.Range("I4").Value = 0.1217029912931362
$ws.Range("J4").Value = 0.1217029912931362
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 398.7037901462904
$ws.Range("R4").Value = 3588.334111316614
$ws.Range("S4").Value = 0.03905987157842264
$ws.Range("T4").Value = 0.03905987157842264

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.588675333333333
$ws.Range("H5").Value = 10.766026
$ws.Range("I5").Value = 0.1217029912931362
$ws.Range("J5").Value = 0.1217029912931362
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 108.1735817272298
$ws.Range("R5").Value = 973.5622355450679
$ws.Range("S5").Value = 0.01059745684607938
$ws.Range("T5").Value = 0.01059745684607938

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.16892433333334
$ws.Range("H6").Value = 54.50677300000001
$ws.Range("I6").Value = 0.6161639698655711
$ws.Range("J6").Value = 0.6161639698655711
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 23.63579766666667
$ws.Range("N6").Value = 70.907393
$ws.Range("O6").Value = 0.06827844587621175
$ws.Range("P6").Value = 0.06827844587621175
$ws.Range("Q6").Value = 429.4370193636433
$ws.Range("R6").Value = 3864.93317427279
$ws.Range("S6").Value = 0.04207071826733817
$ws.Range("T6").Value = 0.04207071826733817

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.16892433333334
$ws.Range("H7").Value = 54.50677300000001
$ws.Range("I7").Value = 0.6161639698655711
$ws.Range("J7").Value = 0.6161639698655711
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 181.2883913333334
$ws.Range("N7").Value = 543.865174
$ws.Range("O7").Value = 0.5237009467675041
$ws.Range("P7").Value = 0.523700946767504
$ws.Range("Q7").Value = 3293.815064647057
$ws.Range("R7").Value = 29644.33558182351
$ws.Range("S7").Value = 0.3226856543826235
$ws.Range("T7").Value = 0.3226856543826234

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.16892433333334
$ws.Range("H8").Value = 54.50677300000001
$ws.Range("I8").Value = 0.6161639698655711
$ws.Range("J8").Value = 0.6161639698655711
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 111.1005463333333
$ws.Range("N8").Value = 333.301639
$ws.Range("O8").Value = 0.3209442197221123
$ws.Range("P8").Value = 0.3209442197221123
$ws.Range("Q8").Value = 2018.577419722328
$ws.Range("R8").Value = 18167.19677750095
$ws.Range("S8").Value = 0.1977542645293848
$ws.Range("T8").Value = 0.1977542645293848

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.16892433333334
$ws.Range("H9").Value = 54.50677300000001
$ws.Range("I9").Value = 0.6161639698655711
$ws.Range("J9").Value = 0.6161639698655711
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.14303933333333
$ws.Range("N9").Value = 90.42911799999999
$ws.Range("O9").Value = 0.08707638763417187
$ws.Range("P9").Value = 0.08707638763417187
$ws.Range("Q9").Value = 547.6666008240238
$ws.Range("R9").Value = 4928.999407416214
$ws.Range("S9").Value = 0.05365333268622467
$ws.Range("T9").Value = 0.05365333268622467

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.554717666666666
$ws.Range("H10").Value = 7.664153
$ws.Range("I10").Value = 0.08663831443731085
$ws.Range("J10").Value = 0.08663831443731086
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.63579766666667
$ws.Range("N10").Value = 70.907393
$ws.Range("O10").Value = 0.06827844587621175
$ws.Range("P10").Value = 0.06827844587621175
$ws.Range("Q10").Value = 60.3827898647921
$ws.Range("R10").Value = 543.4451087831289
$ws.Range("S10").Value = 0.005915529463114143
$ws.Range("T10").Value = 0.005915529463114144

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.554717666666666
$ws.Range("H11").Value = 7.664153
$ws.Range("I11").Value = 0.08663831443731085
$ws.Range("J11").Value = 0.08663831443731086
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 181.2883913333334
$ws.Range("N11").Value = 543.865174
$ws.Range("O11").Value = 0.5237009467675041
$ws.Range("P11").Value = 0.523700946767504
$ws.Range("Q11").Value = 463.1406561008469
$ws.Range("R11").Value = 4168.265904907622
$ws.Range("S11").Value = 0.04537256729716041
$ws.Range("T11").Value = 0.04537256729716041

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.554717666666666
$ws.Range("H12").Value = 7.664153
$ws.Range("I12").Value = 0.08663831443731085
$ws.Range("J12").Value = 0.08663831443731086
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 111.1005463333333
$ws.Range("N12").Value = 333.301639
$ws.Range("O12").Value = 0.3209442197221123
$ws.Range("P12").Value = 0.3209442197221123
$ws.Range("Q12").Value = 283.8305284940852
$ws.Range("R12").Value = 2554.474756446767
$ws.Range("S12").Value = 0.02780606622512174
$ws.Range("T12").Value = 0.02780606622512175

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.554717666666666
$ws.Range("H13").Value = 7.664153
$ws.Range("I13").Value = 0.08663831443731085
$ws.Range("J13").Value = 0.08663831443731086
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 30.14303933333333
$ws.Range("N13").Value = 90.42911799999999
$ws.Range("O13").Value = 0.08707638763417187
$ws.Range("P13").Value = 0.08707638763417187
$ws.Range("Q13").Value = 77.00695511189488
$ws.Range("R13").Value = 693.0625960070539
$ws.Range("S13").Value = 0.007544151451914549
$ws.Range("T13").Value = 0.00754415145191455

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 5.174840666666667
$ws.Range("H14").Value = 15.524522
$ws.Range("I14").Value = 0.1754947244039818
$ws.Range("J14").Value = 0.1754947244039818
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 23.63579766666667
$ws.Range("N14").Value = 70.907393
$ws.Range("O14").Value = 0.06827844587621175
$ws.Range("P14").Value = 0.06827844587621175
$ws.Range("Q14").Value = 122.3114869545718
$ws.Range("R14").Value = 1100.803382591146
$ws.Range("S14").Value = 0.01198250704177797
$ws.Range("T14").Value = 0.01198250704177797

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 5.174840666666667
$ws.Range("H15").Value = 15.524522
$ws.Range("I15").Value = 0.1754947244039818
$ws.Range("J15").Value = 0.1754947244039818
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 181.2883913333334
$ws.Range("N15").Value = 543.865174
$ws.Range("O15").Value = 0.5237009467675041
$ws.Range("P15").Value = 0.523700946767504
$ws.Range("Q15").Value = 938.1385398663144
$ws.Range("R15").Value = 8443.246858796829
$ws.Range("S15").Value = 0.09190675332306747
$ws.Range("T15").Value = 0.09190675332306746

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 5.174840666666667
$ws.Range("H16").Value = 15.524522
$ws.Range("I16").Value = 0.1754947244039818
$ws.Range("J16").Value = 0.1754947244039818
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 111.1005463333333
$ws.Range("N16").Value = 333.301639
$ws.Range("O16").Value = 0.3209442197221123
$ws.Range("P16").Value = 0.3209442197221123
$ws.Range("Q16").Value = 574.9276252546176
$ws.Range("R16").Value = 5174.348627291558
$ws.Range("S16").Value = 0.05632401738918307
$ws.Range("T16").Value = 0.05632401738918307

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 5.174840666666667
$ws.Range("H17").Value = 15.524522
$ws.Range("I17").Value = 0.1754947244039818
$ws.Range("J17").Value = 0.1754947244039818
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 30.14303933333333
$ws.Range("N17").Value = 90.42911799999999
$ws.Range("O17").Value = 0.08707638763417187
$ws.Range("P17").Value = 0.08707638763417187
$ws.Range("Q17").Value = 155.9854257590662
$ws.Range("R17").Value = 1403.868831831596
$ws.Range("S17").Value = 0.01528144664995328
$ws.Range("T17").Value = 0.01528144664995328
